$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 30, shifting existing rows 30-32 down to 31-33
$ws.Rows("30:30").Insert()

# Populate the new row 30 with the new data point.
# Columns A,B,C,E,F,G,H,I,N,O,Q,R are identical to the other "Ramas de apio" rows,
# copy formatting/values from the row below (old row 30, now row 31).
$ws.Range("A31:R31").Copy()
$ws.Range("A30:R30").PasteSpecial(-4104)  # xlPasteAll
$excel.CutCopyMode = 0

$ws.Range("D30").Value = 45194
$ws.Range("J30").Value = 40
$ws.Range("K30").Value = 6000
$ws.Range("L30").Value = 6000
$ws.Range("M30").Value = 6000
$ws.Range("P30").Value = 6000
